$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 609, shifting existing rows 609..664 down to 610..665.
$ws.Rows.Item(609).Insert()

# Populate the newly inserted row 609 with the new record's data.
$ws.Cells.Item(609, 1).Value = 3
$ws.Cells.Item(609, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(609, 3).Value = "Coquimbo"
$ws.Cells.Item(609, 4).Value = 45166
$ws.Cells.Item(609, 5).Value = 5
$ws.Cells.Item(609, 6).Value = 100112031
$ws.Cells.Item(609, 7).Value = "Poroto verde"
$ws.Cells.Item(609, 8).Value = "Sin especificar"
$ws.Cells.Item(609, 9).Value = "Primera"
$ws.Cells.Item(609, 10).Value = 76
$ws.Cells.Item(609, 11).Value = 33000
$ws.Cells.Item(609, 12).Value = 34000
$ws.Cells.Item(609, 13).Value = 33500
$ws.Cells.Item(609, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(609, 15).Value = "Perú"
$ws.Cells.Item(609, 16).Value = 1340
$ws.Cells.Item(609, 17).Value = 25
$ws.Cells.Item(609, 18).Value = "Hortaliza"
